$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 8 formatting: it currently carries an explicit "whole row"
# style (customFormat) left over from how the sheet was produced. The
# other shaded data rows (2, 4, 6) only carry per-cell styles, so bring
# row 8 in line with them by clearing its row-level format and then
# copying the cell-level formatting straight from row 6 (identical
# pattern: Model/Parameters/BIC.../Rsq... shading) back onto row 8.
$ws.Rows.Item(8).ClearFormats()
$srcFormatRow = $ws.Range("A6:F6")
$dstFormatRow = $ws.Range("A8:F8")
$srcFormatRow.Copy()
$dstFormatRow.PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Rename the existing (data) sheet to match the new naming scheme
# used for this export, then insert a new, blank sheet ahead of it for
# the current modeling run.
$ws.Name = "091123"

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "011424"
